$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.826.92'
$ws.Range('E2').Value = '  +0.43%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.753.93'
$ws.Range('E3').Value = '  +0.46%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.01'
$ws.Range('E5').Value = '  -0.70%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  +0.02%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5090'
$ws.Range('E7').Value = '  +3.09%  '

$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2702'
$ws.Range('E8').Value = '  +9.22%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06202'
$ws.Range('E9').Value = '  +3.82%  '

$ws.Range('B10').Value = 'WrappedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.744.96'
$ws.Range('E10').Value = '  -0.04%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06922'
$ws.Range('E11').Value = '  +1.95%  '

$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.58'
$ws.Range('E12').Value = '  +5.05%  '

$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6281'
$ws.Range('E13').Value = '  +7.74%  '

$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '78.52'
$ws.Range('E14').Value = '  +1.71%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.496'
$ws.Range('E15').Value = '  +0.59%  '

$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  +0.06%  '

$ws.Range('B17').Value = 'Dai'
$ws.Range('C17').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9997'
$ws.Range('E17').Value = '  -0.03%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.844.13'
$ws.Range('E18').Value = '  +0.33%  '

$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.72'
$ws.Range('E19').Value = '  +1.46%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006731'
$ws.Range('E20').Value = '  +3.39%  '

$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.970.35'
$ws.Range('E21').Value = '  +0.13%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.076'
$ws.Range('E22').Value = '  +2.51%  '

$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.263'
$ws.Range('E23').Value = '  +4.43%  '

$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.179'
$ws.Range('E24').Value = '  +2.99%  '

$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.63'
$ws.Range('E25').Value = '  +0.35%  '

$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.35'
$ws.Range('E26').Value = '  +5.43%  '

$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.460'
$ws.Range('E27').Value = '  -2.05%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.794'
$ws.Range('E28').Value = '  -2.26%  '

$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '102.86'
$ws.Range('E29').Value = '  +1.88%  '

$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08271'
$ws.Range('E30').Value = '  +2.03%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.730'
$ws.Range('E31').Value = '  -2.01%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.429'
$ws.Range('E32').Value = '  +2.33%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04416'
$ws.Range('E33').Value = '  +0.01%  '

$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9989'
$ws.Range('E34').Value = '  -0.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.645'
$ws.Range('E35').Value = '  -0.11%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.005'
$ws.Range('E36').Value = '  -1.21%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6059'
$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('E38').Value = '  -0.03%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.963'
$ws.Range('E39').Value = '  -4.65%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01563'
$ws.Range('E40').Value = '  +4.27%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.000'
$ws.Range('E41').Value = '  +0.03%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '102.30'
$ws.Range('E42').Value = '  -1.21%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3866'
$ws.Range('E43').Value = '  +2.55%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7510'
$ws.Range('E44').Value = '  -3.51%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.884'
$ws.Range('E45').Value = '  -5.84%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05508'
$ws.Range('E46').Value = '  +7.48%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1095'
$ws.Range('E47').Value = '  +1.25%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.971'
$ws.Range('E48').Value = '  +0.16%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.27'
$ws.Range('E49').Value = '  -0.17%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.94'
$ws.Range('E50').Value = '  +0.63%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.002'
$ws.Range('E51').Value = '  +0.48%  '
